# Add SkillType enum block to Sheet1, mirroring the existing AttributeType / SkillId blocks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 14: new enum declaration row (like rows 4 and 10)
$ws.Range("B14").Value = "SkillType"
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = $true

$ws.Range("G14").Value = "ACTIVE"
$ws.Range("I14").Value = 0

# Row 15: second enum value row
$ws.Range("G15").Value = "PASSIVE"
$ws.Range("I15").Value = 1

# Update the selection to match the final cursor position from the diff
$ws.Range("I15").Select()
